$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 74
$ws1.Range("F3").Value = 299
$ws1.Range("F4").Value = 4357
$ws1.Range("F6").Value = 463

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 74
$ws4.Range("F3").Value = 299
$ws4.Range("F4").Value = 4357
$ws4.Range("F8").Value = 463
